# Trade #98 closed at 2026-02-17 09:17:36 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# --- Summary sheet: totals roll forward ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = 98      # Total Trades
$wsSummary.Range("B9").Value = 41.84   # Win Rate %

# --- Strategy Status sheet: MarketMaking row (row 4) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D4").Value = 98       # Trades
$wsStatus.Range("G4").Value = 41.84    # Win Rate %

# --- New trade row appended as row 99 on both trade-log sheets ---
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 99

    $ws.Cells.Item($row, 1).Value = 98                                  # A - Trade #

    # B - Date: force text storage (matches source data, which stores
    # dates as plain strings, not Excel date serials) by flipping the
    # cell to Text format before the write, then stripping the format
    # back off so no stray number-format sticks to the cell.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).ClearFormats()

    $ws.Cells.Item($row, 3).Value = "09:17:30"                          # C - Time
    $ws.Cells.Item($row, 4).Value = "MarketMaking"                      # D - Strategy
    $ws.Cells.Item($row, 5).Value = "UP"                                # E - Side
    $ws.Cells.Item($row, 6).Value = 0.14                                # F - Entry Price
    $ws.Cells.Item($row, 7).Value = 0.140114                            # G - Exit Price
    $ws.Cells.Item($row, 8).Value = "CLOSED"                            # H - Status
    $ws.Cells.Item($row, 9).Value = 0.08110000000000001                 # I - P&L %
    $ws.Cells.Item($row, 10).Value = 0                                  # J - P&L $
    $ws.Cells.Item($row, 11).Value = 100.11                             # K - Capital After
    $ws.Cells.Item($row, 12).Value = 0                                  # L - Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                                  # M - Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6                                # N - Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps" # O - Entry Reason
    $ws.Cells.Item($row, 16).Value = "early_exit"                       # P - Exit Reason
    $ws.Cells.Item($row, 17).Value = 0.14                               # Q - Duration (min)
}
